$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 487.05
$ws.Range("I15").Value = 487.05
$ws.Range("K15").Value = 1461.15
$ws.Range("M15").Value = -1292.15
$ws.Range("H32").Value = 1778.4
$ws.Range("J32").Value = 2001
$ws.Range("L32").Value = 2001
$ws.Range("N32").Value = -2653
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H111").Value = 2560.5454
$ws.Range("I111").Value = 2518.1667
$ws.Range("J111").Value = 2611.4
$ws.Range("K111").Value = 7554.500100000001
$ws.Range("L111").Value = 7834.200000000001
$ws.Range("M111").Value = -4487.500100000001
$ws.Range("N111").Value = -13968.2
$ws.Range("H127").Value = 697.5
$ws.Range("I127").Value = 697.5
$ws.Range("K127").Value = 2092.5
$ws.Range("M127").Value = 2867.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 2774.5
$ws.Range("I11").Value = 250
$ws.Range("J11").Value = 3279.4
$ws.Range("K11").Value = 250
$ws.Range("L11").Value = 3279.4
$ws.Range("M11").Value = -106
$ws.Range("N11").Value = -3567.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2329.9285
$ws.Range("I86").Value = 2078.4443
$ws.Range("J86").Value = 2782.6
$ws.Range("K86").Value = 2078.4443
$ws.Range("L86").Value = 2782.6
$ws.Range("M86").Value = -955.4443000000001
$ws.Range("N86").Value = -5028.6
$ws.Range("H89").Value = 2329.9285
$ws.Range("I89").Value = 2078.4443
$ws.Range("J89").Value = 2782.6
$ws.Range("K89").Value = 10392.2215
$ws.Range("L89").Value = 13913
$ws.Range("M89").Value = -4776.2215
$ws.Range("N89").Value = -25145
$ws.Range("H127").Value = 89390
$ws.Range("J127").Value = 89390
$ws.Range("L127").Value = 89390
$ws.Range("N127").Value = -99310

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1681.8125
$ws.Range("I7").Value = 1512.5714
$ws.Range("K7").Value = 1512.5714
$ws.Range("M7").Value = -1399.5714
$ws.Range("H99").Value = 7251.923
$ws.Range("I99").Value = 6864.143
$ws.Range("J99").Value = 7704.3335
$ws.Range("K99").Value = 6864.143
$ws.Range("L99").Value = 7704.3335
$ws.Range("M99").Value = -5366.143
$ws.Range("N99").Value = -10700.3335
$ws.Range("H107").Value = 634.8182
$ws.Range("I107").Value = 220.55556
$ws.Range("J107").Value = 2499
$ws.Range("K107").Value = 220.55556
$ws.Range("L107").Value = 2499
$ws.Range("M107").Value = 1699.44444
$ws.Range("N107").Value = -6339
$ws.Range("H126").Value = 7251.923
$ws.Range("I126").Value = 6864.143
$ws.Range("J126").Value = 7704.3335
$ws.Range("K126").Value = 20592.429
$ws.Range("L126").Value = 23113.0005
$ws.Range("M126").Value = -18122.429
$ws.Range("N126").Value = -28053.0005
$ws.Range("H132").Value = 6337.409
$ws.Range("I132").Value = 2207.9167
$ws.Range("K132").Value = 6623.750100000001
$ws.Range("M132").Value = -4093.750100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1466.1111
$ws.Range("I18").Value = 798.75
$ws.Range("K18").Value = 2396.25
$ws.Range("M18").Value = -2227.25
$ws.Range("H75").Value = 607.8333
$ws.Range("I75").Value = 765.6667
$ws.Range("J75").Value = 450
$ws.Range("K75").Value = 2297.0001
$ws.Range("L75").Value = 1350
$ws.Range("M75").Value = -1299.0001
$ws.Range("N75").Value = -3346
$ws.Range("H78").Value = 607.8333
$ws.Range("I78").Value = 765.6667
$ws.Range("J78").Value = 450
$ws.Range("K78").Value = 6891.0003
$ws.Range("L78").Value = 4050
$ws.Range("M78").Value = -1899.0003
$ws.Range("N78").Value = -14034
$ws.Range("H98").Value = 624.3333
$ws.Range("I98").Value = 552.6667
$ws.Range("J98").Value = 696
$ws.Range("K98").Value = 1658.0001
$ws.Range("L98").Value = 2088
$ws.Range("M98").Value = -160.0001
$ws.Range("N98").Value = -5084
$ws.Range("H131").Value = 1733.341
$ws.Range("J131").Value = 1745.5853
$ws.Range("L131").Value = 5236.7559
$ws.Range("N131").Value = -15316.7559
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 15999
$ws.Range("J92").Value = 15999
$ws.Range("L92").Value = 15999
$ws.Range("N92").Value = -19743
$ws.Range("H95").Value = 14999.5
$ws.Range("J95").Value = 14999.5
$ws.Range("L95").Value = 14999.5
$ws.Range("N95").Value = -20491.5
$ws.Range("H102").Value = 711
$ws.Range("I102").Value = 711
$ws.Range("K102").Value = 711
$ws.Range("M102").Value = 911
$ws.Range("H126").Value = 500001440
$ws.Range("I126").Value = 500001440
$ws.Range("K126").Value = 1500004320
$ws.Range("M126").Value = -1500001850
$ws.Range("H132").Value = 2695.5417
$ws.Range("I132").Value = 2223.5334
$ws.Range("K132").Value = 6670.600199999999
$ws.Range("M132").Value = -4140.600199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8131
$ws.Range("I40").Value = 5824.5
$ws.Range("K40").Value = 5824.5
$ws.Range("M40").Value = -5688.5
$ws.Range("H82").Value = 3217
$ws.Range("I82").Value = 2088.5715
$ws.Range("J82").Value = 4533.5
$ws.Range("K82").Value = 2088.5715
$ws.Range("L82").Value = 4533.5
$ws.Range("M82").Value = -1727.5715
$ws.Range("N82").Value = -5255.5
$ws.Range("H85").Value = 3217
$ws.Range("I85").Value = 2088.5715
$ws.Range("J85").Value = 4533.5
$ws.Range("K85").Value = 2088.5715
$ws.Range("L85").Value = 4533.5
$ws.Range("M85").Value = -840.5715
$ws.Range("N85").Value = -7029.5
$ws.Range("H93").Value = 848.0833
$ws.Range("I93").Value = 729.2
$ws.Range("J93").Value = 1442.5
$ws.Range("K93").Value = 729.2
$ws.Range("L93").Value = 1442.5
$ws.Range("M93").Value = 518.8
$ws.Range("N93").Value = -3938.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 619.5
$ws.Range("I107").Value = 572.75
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 1718.25
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = 201.75
$ws.Range("N107").Value = -6540
$ws.Range("H132").Value = 3819.2222
$ws.Range("I132").Value = 3819.2222
$ws.Range("K132").Value = 11457.6666
$ws.Range("M132").Value = -8927.6666
